$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C, shifting the existing
# column C (English translation) data into column D.
$ws.Columns("C:C").Insert()
